$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New diary entry for 30 loka (Oct 30), rows added under the existing table.
$ws.Range("A19").Value = "30 loka"

$ws.Range("B19").Value = "9.15-11.15"
$ws.Range("B19").NumberFormat = "h:mm"

$ws.Range("C19").Value = "Initial draw ongelman selvittelyä fireworkscenestä, oppikirjasta 104-"
$ws.Range("C19").WrapText = $true

$ws.Range("G19").Value = 2

# Match the row height used by the other wrapped-text rows in the log.
$ws.Rows.Item(19).RowHeight = 43.5

# Move the view so the newly added row is visible and select it, like the
# author left the sheet after typing the entry.
try {
    $excel.ActiveWindow.ScrollRow = 12
} catch {}
$ws.Range("H19").Select()
